$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new worksheet "AddPriceAgrmnt_FavFolder" after the last
#    existing sheet (AddPriceAgrmnt_LocalCatalog).
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFav = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsFav.Name = "AddPriceAgrmnt_FavFolder"

# ------------------------------------------------------------------
# 2. Populate "AddPriceAgrmnt_LocalCatalog" (3rd sheet) with the new
#    ItemType / ItemNumber columns.
# ------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("AddPriceAgrmnt_LocalCatalog")
$wsLocal.Range("A1").Value = "Role"
$wsLocal.Range("B1").Value = "Location                      "
$wsLocal.Range("C1").Value = "ItemType "
$wsLocal.Range("D1").Value = "ItemNumber"
$wsLocal.Range("A2").Value = "REQUESTOR"
$wsLocal.Range("B2").Value = "XEEVA -MJ"
$wsLocal.Range("C2").Value = "local"
$wsLocal.Range("D2").Value = "156001-00013"

# Match the authored column widths (best-fit) as closely as possible.
$wsLocal.Columns.Item(1).ColumnWidth = 10.6
$wsLocal.Columns.Item(2).ColumnWidth = 8.76
$wsLocal.Columns.Item(3).ColumnWidth = 8.93
$wsLocal.Columns.Item(4).ColumnWidth = 11.6

# ------------------------------------------------------------------
# 3. Populate the new "AddPriceAgrmnt_FavFolder" sheet with the
#    favourite-folder sample row.
# ------------------------------------------------------------------
$wsFav.Range("A1").Value = "Role"
$wsFav.Range("B1").Value = "Location                      "
$wsFav.Range("A2").Value = "REQUESTOR"
$wsFav.Range("B2").Value = "XEEVA -MJ"
$wsFav.Columns.Item(3).ColumnWidth = 9.93

# ------------------------------------------------------------------
# 4. Selections / active cells to mirror the authored workbook.
# ------------------------------------------------------------------

# AddPriceAgrmnt_RecentOrder (2nd sheet) - selection becomes A1:B2,
# tab no longer selected since a later sheet becomes active.
$wsRecent = $wb.Worksheets.Item("AddPriceAgrmnt_RecentOrder")
$wsRecent.Range("A1:B2").Select() | Out-Null

# AddPriceAgrmnt_LocalCatalog - selection becomes A1:C2
$wsLocal.Range("A1:C2").Select() | Out-Null

# AddPriceAgrmnt_FavFolder - becomes the active sheet with B7 selected
$wsFav.Activate() | Out-Null
$wsFav.Range("B7").Select() | Out-Null
